$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.156959333333335
$ws.Range("H2").Value = 27.470878
$ws.Range("I2").Value = 0.969469463764299
$ws.Range("J2").Value = 0.9694694637642989
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.558821666666667
$ws.Range("N2").Value = 7.676465
$ws.Range("O2").Value = 0.2156728774407755
$ws.Range("P2").Value = 0.2156728774407755
$ws.Range("Q2").Value = 23.43102594291889
$ws.Range("R2").Value = 210.87923348627
$ws.Range("S2").Value = 0.2090882688410119
$ws.Range("T2").Value = 0.2090882688410119

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.156959333333335
$ws.Range("H3").Value = 27.470878
$ws.Range("I3").Value = 0.969469463764299
$ws.Range("J3").Value = 0.9694694637642989
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.453984666666667
$ws.Range("N3").Value = 19.361954
$ws.Range("O3").Value = 0.5439806384912759
$ws.Range("P3").Value = 0.5439806384912759
$ws.Range("Q3").Value = 59.09887513062357
$ws.Range("R3").Value = 531.889876175612
$ws.Range("S3").Value = 0.5273726178962982
$ws.Range("T3").Value = 0.5273726178962982

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd8"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.156959333333335
$ws.Range("H4").Value = 27.470878
$ws.Range("I4").Value = 0.969469463764299
$ws.Range("J4").Value = 0.9694694637642989
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.851558333333334
$ws.Range("N4").Value = 8.554675000000001
$ws.Range("O4").Value = 0.2403464840679487
$ws.Range("P4").Value = 0.2403464840679487
$ws.Range("Q4").Value = 26.11160369496112
$ws.Range("R4").Value = 235.0044332546501
$ws.Range("S4").Value = 0.2330085770269889
$ws.Range("T4").Value = 0.2330085770269889

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd8"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.288371
$ws.Range("H5").Value = 0.865113
$ws.Range("I5").Value = 0.03053053623570109
$ws.Range("J5").Value = 0.03053053623570109
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.558821666666667
$ws.Range("N5").Value = 7.676465
$ws.Range("O5").Value = 0.2156728774407755
$ws.Range("P5").Value = 0.2156728774407755
$ws.Range("Q5").Value = 0.7378899628383333
$ws.Range("R5").Value = 6.641009665545001
$ws.Range("S5").Value = 0.006584608599763515
$ws.Range("T5").Value = 0.006584608599763515

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd8"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.288371
$ws.Range("H6").Value = 0.865113
$ws.Range("I6").Value = 0.03053053623570109
$ws.Range("J6").Value = 0.03053053623570109
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.453984666666667
$ws.Range("N6").Value = 19.361954
$ws.Range("O6").Value = 0.5439806384912759
$ws.Range("P6").Value = 0.5439806384912759
$ws.Range("Q6").Value = 1.861142012311333
$ws.Range("R6").Value = 16.750278110802
$ws.Range("S6").Value = 0.01660802059497771
$ws.Range("T6").Value = 0.01660802059497771

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd8"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.288371
$ws.Range("H7").Value = 0.865113
$ws.Range("I7").Value = 0.03053053623570109
$ws.Range("J7").Value = 0.03053053623570109
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.851558333333334
$ws.Range("N7").Value = 8.554675000000001
$ws.Range("O7").Value = 0.2403464840679487
$ws.Range("P7").Value = 0.2403464840679487
$ws.Range("Q7").Value = 0.8223067281416667
$ws.Range("R7").Value = 7.400760553275001
$ws.Range("S7").Value = 0.007337907040959862
$ws.Range("T7").Value = 0.007337907040959863
